$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 0.04856454963496049
$ws.Range("D2").Value = 0.1511996330835572
$ws.Range("E2").Value = 0.1364315755092314
$ws.Range("F2").Value = 1.544159145771673
$ws.Range("G2").Value = 0.00246662241494966
$ws.Range("J2").Value = 0.1576926278446074
$ws.Range("K2").Value = 1.116889796214252
$ws.Range("M2").Value = 0.3924742335102991
$ws.Range("N2").Value = 1.413550467732577
$ws.Range("O2").Value = 3.855768653390697

# Row 3
$ws.Range("C3").Value = 0.04316241441010504
$ws.Range("D3").Value = 0.1496243054887003
$ws.Range("E3").Value = 0.1362465658777374
$ws.Range("F3").Value = 1.544578384392402
$ws.Range("G3").Value = 0.002469760370666819
$ws.Range("J3").Value = 0.1583497748499916
$ws.Range("K3").Value = 1.008479844026823
$ws.Range("M3").Value = 0.3698125674116568
$ws.Range("N3").Value = 1.425553101108861
$ws.Range("O3").Value = 3.862406596950706

# Row 4
$ws.Range("C4").Value = 0.03986008640386274
$ws.Range("D4").Value = 0.1487056145324246
$ws.Range("E4").Value = 0.1361899399102064
$ws.Range("F4").Value = 1.545732288601208
$ws.Range("G4").Value = 0.002471790567627382
$ws.Range("J4").Value = 0.1588294235726977
$ws.Range("K4").Value = 0.9420370253374131
$ws.Range("M4").Value = 0.3560254197317505
$ws.Range("N4").Value = 1.433452495177939
$ws.Range("O4").Value = 3.868867740159089

# Row 5
$ws.Range("C5").Value = 0.03851801591821413
$ws.Range("D5").Value = 0.1483435094232419
$ws.Range("E5").Value = 0.1361812106161775
$ws.Range("F5").Value = 1.54642783116266
$ws.Range("G5").Value = 0.002472643989194912
$ws.Range("J5").Value = 0.1590440345663069
$ws.Range("K5").Value = 0.9149929189295278
$ws.Range("M5").Value = 0.350439334071865
$ws.Range("N5").Value = 1.436804868281136
$ws.Range("O5").Value = 3.872099937677945

# Row 6
$ws.Range("C6").Value = 0.03829538638122187
$ws.Range("D6").Value = 0.1482841248221405
$ws.Range("E6").Value = 0.136180628207283
$ws.Range("F6").Value = 1.546556930834981
$ws.Range("G6").Value = 0.002472787278096024
$ws.Range("J6").Value = 0.1590808271830007
$ws.Range("K6").Value = 0.9105042317784182
$ws.Range("M6").Value = 0.3495137278925995
$ws.Range("N6").Value = 1.437369581707863
$ws.Range("O6").Value = 3.872672816008361

# Row 7
$ws.Range("C7").Value = 0.03984197199430639
$ws.Range("D7").Value = 0.1487006813034526
$ws.Range("E7").Value = 0.1361897640690941
$ws.Range("F7").Value = 1.54574075678228
$ws.Range("G7").Value = 0.00247180197139317
$ws.Range("J7").Value = 0.1588322403530498
$ws.Range("K7").Value = 0.9416721681508307
$ws.Range("M7").Value = 0.3559499526894143
$ws.Range("N7").Value = 1.433497166568387
$ws.Range("O7").Value = 3.868908905258593

# Row 8
$ws.Range("C8").Value = 0.04669885002060425
$ws.Range("D8").Value = 0.1506464132096568
$ws.Range("E8").Value = 0.1363559718356342
$ws.Range("F8").Value = 1.544117560426727
$ws.Range("G8").Value = 0.00246768295261603
$ws.Range("J8").Value = 0.1579034051910178
$ws.Range("K8").Value = 1.07948565033206
$ws.Range("M8").Value = 0.3846342623225993
$ws.Range("N8").Value = 1.417579090365535
$ws.Range("O8").Value = 3.857561949122442

# Row 9
$ws.Range("C9").Value = 0.06026252603402327
$ws.Range("D9").Value = 0.1548451703329476
$ws.Range("E9").Value = 0.1371332032796282
$ws.Range("F9").Value = 1.548054924552886
$ws.Range("G9").Value = 0.00246042301832032
$ws.Range("J9").Value = 0.1566863644232903
$ws.Range("K9").Value = 1.350653619515867
$ws.Range("M9").Value = 0.441883296146969
$ws.Range("N9").Value = 1.39056258646017
$ws.Range("O9").Value = 3.85426979683038

# Row 10
$ws.Range("C10").Value = 0.07030233460730528
$ws.Range("D10").Value = 0.1581611336089708
$ws.Range("E10").Value = 0.1379786067399493
$ws.Range("F10").Value = 1.555301564897718
$ws.Range("G10").Value = 0.002455582420235667
$ws.Range("J10").Value = 0.1561609996018731
$ws.Range("K10").Value = 1.550396937049186
$ws.Range("M10").Value = 0.4845442672287348
$ws.Range("N10").Value = 1.373267517510442
$ws.Range("O10").Value = 3.863460161107781

# Row 11
$ws.Range("C11").Value = 0.07488660156008109
$ws.Range("D11").Value = 0.1597193494405929
$ws.Range("E11").Value = 0.1384226381872367
$ws.Range("F11").Value = 1.559546679025416
$ws.Range("G11").Value = 0.002453486339505449
$ws.Range("J11").Value = 0.1560021644736764
$ws.Range("K11").Value = 1.641369773200267
$ws.Range("M11").Value = 0.5040805396450949
$ws.Range("N11").Value = 1.365953015439302
$ws.Range("O11").Value = 3.87017340342237

# Row 12
$ws.Range("C12").Value = 0.07662503775046048
$ws.Range("D12").Value = 0.1603165115604526
$ws.Range("E12").Value = 0.1385993144778155
$ws.Range("F12").Value = 1.561290790756402
$ws.Range("G12").Value = 0.002452707759632196
$ws.Range("J12").Value = 0.1559535490889701
$ws.Range("K12").Value = 1.675833339128133
$ws.Range("M12").Value = 0.5114968142977077
$ws.Range("N12").Value = 1.363262674442886
$ws.Range("O12").Value = 3.873080523720006

# Row 13
$ws.Range("C13").Value = 0.07625052417468225
$ws.Range("D13").Value = 0.1601875871607206
$ws.Range("E13").Value = 0.1385608849362576
$ws.Range("F13").Value = 1.560909087910346
$ws.Range("G13").Value = 0.002452874767659003
$ws.Range("J13").Value = 0.1559635062996207
$ws.Range("K13").Value = 1.668410388910218
$ws.Range("M13").Value = 0.5098987775212223
$ws.Range("N13").Value = 1.363838552119105
$ws.Range("O13").Value = 3.872438178987352

# Row 14
$ws.Range("C14").Value = 0.07502957409597855
$ws.Range("D14").Value = 0.159768336366767
$ws.Range("E14").Value = 0.1384370025768646
$ws.Range("F14").Value = 1.559687429896741
$ws.Range("G14").Value = 0.002453421981716334
$ws.Range("J14").Value = 0.1559979337043131
$ws.Range("K14").Value = 1.644204836889685
$ws.Range("M14").Value = 0.5046903158565854
$ws.Range("N14").Value = 1.365730086301554
$ws.Range("O14").Value = 3.870405255057847

# Row 15
$ws.Range("C15").Value = 0.07428202948011631
$ws.Range("D15").Value = 0.1595124560423784
$ws.Range("E15").Value = 0.1383622315362771
$ws.Range("F15").Value = 1.558956920910816
$ws.Range("G15").Value = 0.002453759138498163
$ws.Range("J15").Value = 0.1560205234655143
$ws.Range("K15").Value = 1.629380033039752
$ws.Range("M15").Value = 0.5015023579713613
$ws.Range("N15").Value = 1.366899058210699
$ws.Range("O15").Value = 3.869207582667315

# Row 16
$ws.Range("C16").Value = 0.07000308741237404
$ws.Range("D16").Value = 0.1580602971200875
$ws.Range("E16").Value = 0.1379507830430988
$ws.Range("F16").Value = 1.555043239468546
$ws.Range("G16").Value = 0.002455721529523753
$ws.Range("J16").Value = 0.1561729932790783
$ws.Range("K16").Value = 1.544453706121828
$ws.Range("M16").Value = 0.4832701059029674
$ws.Range("N16").Value = 1.373756664351802
$ws.Range("O16").Value = 3.863072462849544

# Row 17
$ws.Range("C17").Value = 0.06738248660241197
$ws.Range("D17").Value = 0.1571821507539681
$ws.Range("E17").Value = 0.1377135872131809
$ws.Range("F17").Value = 1.552885400861399
$ws.Range("G17").Value = 0.002456952472582799
$ws.Range("J17").Value = 0.1562870624559309
$ws.Range("K17").Value = 1.492380908791063
$ws.Range("M17").Value = 0.4721181698864569
$ws.Range("N17").Value = 1.378105222280148
$ws.Range("O17").Value = 3.859957961457695

# Row 18
$ws.Range("C18").Value = 0.06587679578915129
$ws.Range("D18").Value = 0.1566817519335899
$ws.Range("E18").Value = 0.1375827555350178
$ws.Range("F18").Value = 1.551733543629538
$ws.Range("G18").Value = 0.002457670453464643
$ws.Range("J18").Value = 0.1563602163406372
$ws.Range("K18").Value = 1.462440323563669
$ws.Range("M18").Value = 0.4657160873337745
$ws.Range("N18").Value = 1.380658460943103
$ws.Range("O18").Value = 3.858404906029989

# Row 19
$ws.Range("C19").Value = 0.06536727075241799
$ws.Range("D19").Value = 0.1565131323521456
$ws.Range("E19").Value = 0.1375394200343685
$ws.Range("F19").Value = 1.551358872242588
$ws.Range("G19").Value = 0.002457915264928114
$ws.Range("J19").Value = 0.1563862805597935
$ws.Range("K19").Value = 1.452304772168532
$ws.Range("M19").Value = 0.4635505601257464
$ws.Range("N19").Value = 1.381531886735139
$ws.Range("O19").Value = 3.857919976375882

# Row 20
$ws.Range("C20").Value = 0.06766128736050803
$ws.Range("D20").Value = 0.1572751461829967
$ws.Range("E20").Value = 0.1377382579677864
$ws.Range("F20").Value = 1.553105865628922
$ws.Range("G20").Value = 0.002456820404501738
$ws.Range("J20").Value = 0.1562741387487634
$ws.Range("K20").Value = 1.497923090667086
$ws.Range("M20").Value = 0.4733040505272328
$ws.Range("N20").Value = 1.377636922789328
$ws.Range("O20").Value = 3.860264834137382

# Row 21
$ws.Range("C21").Value = 0.07538812948507712
$ws.Range("D21").Value = 0.1598912882584358
$ws.Range("E21").Value = 0.1384731584375984
$ws.Range("F21").Value = 1.560042552356933
$ws.Range("G21").Value = 0.002453260840907692
$ws.Range("J21").Value = 0.1559875085190399
$ws.Range("K21").Value = 1.651314222699853
$ws.Range("M21").Value = 0.5062196728693564
$ws.Range("N21").Value = 1.365172339440633
$ws.Range("O21").Value = 3.87099246284123

# Row 22
$ws.Range("C22").Value = 0.08045249274363186
$ws.Range("D22").Value = 0.1616424439915534
$ws.Range("E22").Value = 0.139003176523584
$ws.Range("F22").Value = 1.565372247326067
$ws.Range("G22").Value = 0.002451022793406382
$ws.Range("J22").Value = 0.1558673993293169
$ws.Range("K22").Value = 1.751645852063518
$ws.Range("M22").Value = 0.5278385244341877
$ws.Range("N22").Value = 1.357489401369641
$ws.Range("O22").Value = 3.880131229602597

# Row 23
$ws.Range("C23").Value = 0.07774822239915125
$ws.Range("D23").Value = 0.1607040539225579
$ws.Range("E23").Value = 0.1387157526306986
$ws.Range("F23").Value = 1.562454781491866
$ws.Range("G23").Value = 0.00245220922277779
$ws.Range("J23").Value = 0.1559253512939804
$ws.Range("K23").Value = 1.698089992150415
$ws.Range("M23").Value = 0.5162904894820173
$ws.Range("N23").Value = 1.361547539923102
$ws.Range("O23").Value = 3.875058756656841

# Row 24
$ws.Range("C24").Value = 0.06753523865164368
$ws.Range("D24").Value = 0.1572330890563336
$ws.Range("E24").Value = 0.1377270870764704
$ws.Range("F24").Value = 1.553005917180968
$ws.Range("G24").Value = 0.002456880080552838
$ws.Range("J24").Value = 0.1562799579649266
$ws.Range("K24").Value = 1.495417480222727
$ws.Range("M24").Value = 0.4727678848534325
$ws.Range("N24").Value = 1.37784847528264
$ws.Range("O24").Value = 3.860125357204453

# Row 25
$ws.Range("C25").Value = 0.05658031382823481
$ws.Range("D25").Value = 0.1536685282590255
$ws.Range("E25").Value = 0.1368747051036685
$ws.Range("F25").Value = 1.546226125759958
$ws.Range("G25").Value = 0.002462300034673941
$ws.Range("J25").Value = 0.1569508666215853
$ws.Range("K25").Value = 1.277201814083924
$ws.Range("M25").Value = 0.441883296146969
$ws.Range("N25").Value = 1.39056258646017
$ws.Range("O25").Value = 3.85426979683038

